$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header additions
$ws.Range("B1").Value = "Untitled 1"
$ws.Range("C1").Value = "Untitled 2"
$ws.Range("D1").Value = "Untitled 3"

# Row 2 data values
$ws.Range("A2").Value = 0.000000
$ws.Range("B2").Value = 0.000007
$ws.Range("C2").Value = 11.562918
$ws.Range("D2").Value = 11.562918

# Match A2's number style on the new row-2 cells (s="1")
$ws.Range("A2").Copy()
$ws.Range("B2:D2").PasteSpecial(-4122)
